$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Obrigatorio) for rows 2-9 from "N" to "S"
$ws.Range("E2:E9").Value = "S"
